$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B, shifting text/tokens/code/sum columns one to the right.
$ws.Range("B1").EntireColumn.Insert()

# The insert copies formatting from the left neighbour (column A), which is the
# opposite of what we want (header styled, data rows plain). Strip that back off
# the data rows, then give the new header cell the same style as the other
# header cells (copy format from C1, the former B1 header).
$ws.Range("B2:B4").ClearFormats()
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "names" column content.
$ws.Range("B1").Value = "names"
$ws.Range("B2").Value = "strvisx"
$ws.Range("B3").Value = "strvis_orig"

# B4's value is the text "4" (not a number) - force text formatting for entry,
# then clear the formatting again so the cell ends up with no explicit style,
# matching the rest of the data column.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "4"
$ws.Range("B4").ClearFormats()
